# Scheduled data refresh: updates currentAveragePrice* / Leve Price / Leve Profit
# columns (H-N) across the per-job Leve tables (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 46
$ws.Range("H46").Value = 8874
$ws.Range("J46").Value = 8874
$ws.Range("L46").Value = 26622
$ws.Range("N46").Value = -26860

# Row 60
$ws.Range("H60").Value = 8874
$ws.Range("J60").Value = 8874
$ws.Range("L60").Value = 26622
$ws.Range("N60").Value = -27590

# Row 96
$ws.Range("H96").Value = 914.2143
$ws.Range("I96").Value = 1576.5714
$ws.Range("J96").Value = 251.85715
$ws.Range("K96").Value = 4729.7142
$ws.Range("L96").Value = 755.5714499999999
$ws.Range("M96").Value = -3356.7142
$ws.Range("N96").Value = -3501.57145

# Row 100
$ws.Range("H100").Value = 1689.2
$ws.Range("I100").Value = 1795.4615
$ws.Range("J100").Value = 998.5
$ws.Range("K100").Value = 1795.4615
$ws.Range("L100").Value = 998.5
$ws.Range("M100").Value = -1254.4615
$ws.Range("N100").Value = -2080.5

# Row 127
$ws.Range("H127").Value = 1847.5
$ws.Range("I127").Value = 1847.5
$ws.Range("J127").Value = 0
$ws.Range("K127").Value = 5542.5
$ws.Range("L127").Value = 0
$ws.Range("M127").Value = -582.5
$ws.Range("N127").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 6408.5537
$ws.Range("I32").Value = 3890.7
$ws.Range("J32").Value = 36622.8
$ws.Range("K32").Value = 3890.7
$ws.Range("L32").Value = 36622.8
$ws.Range("M32").Value = -3603.7
$ws.Range("N32").Value = -37196.8

# Row 38
$ws.Range("H38").Value = 10013.6
$ws.Range("J38").Value = 10021
$ws.Range("L38").Value = 10021
$ws.Range("N38").Value = -10955

# Row 45
$ws.Range("H45").Value = 10182.818
$ws.Range("I45").Value = 11079.7
$ws.Range("J45").Value = 1214
$ws.Range("K45").Value = 11079.7
$ws.Range("L45").Value = 1214
$ws.Range("M45").Value = -10702.7
$ws.Range("N45").Value = -1968

# Row 61
$ws.Range("H61").Value = 3944.2
$ws.Range("I61").Value = 2303.1072
$ws.Range("J61").Value = 5646.074
$ws.Range("K61").Value = 2303.1072
$ws.Range("L61").Value = 5646.074
$ws.Range("M61").Value = -2091.1072
$ws.Range("N61").Value = -6070.074

# Row 88
$ws.Range("H88").Value = 867.26086
$ws.Range("I88").Value = 952.0909
$ws.Range("K88").Value = 952.0909
$ws.Range("M88").Value = -546.0909

# Row 91
$ws.Range("H91").Value = 867.26086
$ws.Range("I91").Value = 952.0909
$ws.Range("K91").Value = 952.0909
$ws.Range("M91").Value = 451.9091

# Row 117
$ws.Range("H117").Value = 37666.332
$ws.Range("J117").Value = 37666.332
$ws.Range("L117").Value = 37666.332
$ws.Range("N117").Value = -46844.332

# Row 132
$ws.Range("H132").Value = 2140.6604
$ws.Range("I132").Value = 1964.9783
$ws.Range("K132").Value = 5894.9349
$ws.Range("M132").Value = -3364.9349

# Row 136
$ws.Range("H136").Value = 3944.2
$ws.Range("I136").Value = 2303.1072
$ws.Range("J136").Value = 5646.074
$ws.Range("K136").Value = 6909.321599999999
$ws.Range("L136").Value = 16938.222
$ws.Range("M136").Value = -4359.321599999999
$ws.Range("N136").Value = -22038.222

$ws = $wb.Worksheets.Item("BSM")
# Row 19
$ws.Range("H19").Value = 20004.5
$ws.Range("I19").Value = 19999
$ws.Range("J19").Value = 20010
$ws.Range("K19").Value = 19999
$ws.Range("L19").Value = 20010
$ws.Range("M19").Value = -19826
$ws.Range("N19").Value = -20356

# Row 20
$ws.Range("H20").Value = 4622.0347
$ws.Range("I20").Value = 3687.1428
$ws.Range("K20").Value = 3687.1428
$ws.Range("M20").Value = -3440.1428

# Row 42
$ws.Range("H42").Value = 389684
$ws.Range("J42").Value = 389684
$ws.Range("L42").Value = 389684
$ws.Range("N42").Value = -390340

# Row 94
$ws.Range("H94").Value = 2151.6667
$ws.Range("I94").Value = 1294.3572
$ws.Range("K94").Value = 1294.3572
$ws.Range("M94").Value = -843.3571999999999

# Row 99
$ws.Range("H99").Value = 7750.9
$ws.Range("I99").Value = 4249.5
$ws.Range("J99").Value = 10085.167
$ws.Range("K99").Value = 4249.5
$ws.Range("L99").Value = 10085.167
$ws.Range("M99").Value = -2751.5
$ws.Range("N99").Value = -13081.167

# Row 105
$ws.Range("H105").Value = 2685.8462
$ws.Range("I105").Value = 2888.9092
$ws.Range("J105").Value = 1569
$ws.Range("K105").Value = 2888.9092
$ws.Range("L105").Value = 1569
$ws.Range("M105").Value = -1141.9092
$ws.Range("N105").Value = -5063

# Row 134
$ws.Range("H134").Value = 1534.4584
$ws.Range("I134").Value = 1497.9011
$ws.Range("J134").Value = 2199.8
$ws.Range("K134").Value = 4493.7033
$ws.Range("L134").Value = 6599.400000000001
$ws.Range("M134").Value = -1958.7033
$ws.Range("N134").Value = -11669.4

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 40954.08
$ws.Range("I31").Value = 49739.855
$ws.Range("J31").Value = 4053.8
$ws.Range("K31").Value = 49739.855
$ws.Range("L31").Value = 4053.8
$ws.Range("M31").Value = -49444.855
$ws.Range("N31").Value = -4643.8

# Row 34
$ws.Range("H34").Value = 40954.08
$ws.Range("I34").Value = 49739.855
$ws.Range("J34").Value = 4053.8
$ws.Range("K34").Value = 49739.855
$ws.Range("L34").Value = 4053.8
$ws.Range("M34").Value = -49537.855
$ws.Range("N34").Value = -4457.8

# Row 80
$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()

# Row 83
$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()

# Row 132
$ws.Range("H132").Value = 3655.6155
$ws.Range("I132").Value = 3454.3684
$ws.Range("K132").Value = 10363.1052
$ws.Range("M132").Value = -7833.1052

# Row 134
$ws.Range("H134").Value = 15219.893
$ws.Range("I134").Value = 6006.0386
$ws.Range("K134").Value = 18018.1158
$ws.Range("M134").Value = -15483.1158

# Row 141
$ws.Range("H141").Value = 200476.33
$ws.Range("I141").Value = 70000
$ws.Range("J141").Value = 216785.88
$ws.Range("K141").Value = 70000
$ws.Range("L141").Value = 216785.88
$ws.Range("M141").Value = -64820
$ws.Range("N141").Value = -227145.88

$ws = $wb.Worksheets.Item("CUL")
# Row 12
$ws.Range("H12").Value = 226.33333
$ws.Range("I12").Value = 136.25
$ws.Range("J12").Value = 252.07143
$ws.Range("K12").Value = 408.75
$ws.Range("L12").Value = 756.21429
$ws.Range("M12").Value = -235.75
$ws.Range("N12").Value = -1102.21429

# Row 129
$ws.Range("H129").Value = 9901650
$ws.Range("J129").Value = 2021.5
$ws.Range("L129").Value = 6064.5
$ws.Range("N129").Value = -16064.5

# Row 136
$ws.Range("H136").Value = 1263124.9
$ws.Range("I136").Value = 9999999
$ws.Range("K136").Value = 29999997
$ws.Range("M136").Value = -29994897

# Row 137
$ws.Range("H137").Value = 2901.4167
$ws.Range("I137").Value = 1285.8
$ws.Range("J137").Value = 4055.4285
$ws.Range("K137").Value = 3857.4
$ws.Range("L137").Value = 12166.2855
$ws.Range("M137").Value = 1242.6
$ws.Range("N137").Value = -22366.2855

# Row 140
$ws.Range("H140").Value = 3142.923
$ws.Range("I140").Value = 2441.6365
$ws.Range("K140").Value = 7324.9095
$ws.Range("M140").Value = -2144.9095

$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Range("H2").Value = 143.6
$ws.Range("J2").Value = 310.25
$ws.Range("L2").Value = 310.25
$ws.Range("N2").Value = -536.25

# Row 49
$ws.Range("H49").Value = 0
$ws.Range("J49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("N49").ClearContents()

# Row 102
$ws.Range("H102").Value = 111112300
$ws.Range("I102").Value = 958.5714
$ws.Range("J102").Value = 500002020
$ws.Range("K102").Value = 958.5714
$ws.Range("L102").Value = 500002020
$ws.Range("M102").Value = 663.4286
$ws.Range("N102").Value = -500005264

# Row 122
$ws.Range("H122").Value = 3249.2222
$ws.Range("I122").Value = 2797.5667
$ws.Range("J122").Value = 5507.5
$ws.Range("K122").Value = 8392.7001
$ws.Range("L122").Value = 16522.5
$ws.Range("M122").Value = -5942.7001
$ws.Range("N122").Value = -21422.5

# Row 132
$ws.Range("H132").Value = 3139.0454
$ws.Range("I132").Value = 2612.4546
$ws.Range("J132").Value = 4718.8184
$ws.Range("K132").Value = 7837.3638
$ws.Range("L132").Value = 14156.4552
$ws.Range("M132").Value = -5307.3638
$ws.Range("N132").Value = -19216.4552

# Row 134
$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 3433.6155
$ws.Range("I40").Value = 2535.2273
$ws.Range("K40").Value = 2535.2273
$ws.Range("M40").Value = -2399.2273

# Row 93
$ws.Range("H93").Value = 1659.2727
$ws.Range("I93").Value = 1625.2
$ws.Range("K93").Value = 1625.2
$ws.Range("M93").Value = -377.2

# Row 112
$ws.Range("H112").Value = 69990
$ws.Range("J112").Value = 69990
$ws.Range("L112").Value = 69990
$ws.Range("N112").Value = -72944

$ws = $wb.Worksheets.Item("WVR")
# Row 100
$ws.Range("H100").Value = 679.89655
$ws.Range("I100").Value = 686.6
$ws.Range("J100").Value = 665
$ws.Range("K100").Value = 1373.2
$ws.Range("L100").Value = 1330
$ws.Range("M100").Value = -832.2
$ws.Range("N100").Value = -2412

# Row 122
$ws.Range("H122").Value = 1661.7142
$ws.Range("I122").Value = 1558.8077
$ws.Range("K122").Value = 4676.4231
$ws.Range("M122").Value = -2226.4231

# Row 132
$ws.Range("H132").Value = 1770.5454
$ws.Range("I132").Value = 1691.2258
$ws.Range("K132").Value = 5073.6774
$ws.Range("M132").Value = -2543.6774

# Row 141
$ws.Range("H141").Value = 63375
$ws.Range("I141").Value = 80000
$ws.Range("J141").Value = 61000
$ws.Range("K141").Value = 80000
$ws.Range("L141").Value = 61000
$ws.Range("M141").Value = -74820
$ws.Range("N141").Value = -71360
